# Generate Report for Handback
#
# The da04f27e-850a-40af-8ab7-b2b0fd945ad0.md file has completed its
# handback cycle: it moves from "Ready for handoff" to
# "Handed back: in sync with en-US" on the Overview sheet and on both
# locale sheets, and the per-locale "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" columns (row 6) get
# populated for zh-cn and de-de.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$statusDone  = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status columns (E6, F6)
# for the da04f27e row.
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E6").Value2 = $statusDone
$overview.Range("F6").Value2 = $statusDone

# ---------------------------------------------------------------
# zh-cn sheet: row 6 is the da04f27e entry.
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C6").Value2 = $statusDone
$zh.Range("J6").Value2 = "da04f27e-850a-40af-8ab7-b2b0fd945ad0.f92697e7084d2cb9073de0ecd739a6b2c473ccca.zh-cn.xlf"
$zh.Range("K6").Value2 = "2016-09-06 04:49:52"
$zh.Hyperlinks.Add($zh.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/77afc0bdbe9b1f438cf63b54a9c0d3c201ca99c4/e2e/da04f27e-850a-40af-8ab7-b2b0fd945ad0.md", "", "", "da04f27e-850a-40af-8ab7-b2b0fd945ad0.md")

# ---------------------------------------------------------------
# de-de sheet: row 6 is the da04f27e entry.
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C6").Value2 = $statusDone
$de.Range("J6").Value2 = "da04f27e-850a-40af-8ab7-b2b0fd945ad0.f92697e7084d2cb9073de0ecd739a6b2c473ccca.de-de.xlf"
$de.Range("K6").Value2 = "2016-09-06 04:50:19"
$de.Hyperlinks.Add($de.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f503da96d4f1e4156996bca99b43870f8091ef40/e2e/da04f27e-850a-40af-8ab7-b2b0fd945ad0.md", "", "", "da04f27e-850a-40af-8ab7-b2b0fd945ad0.md")
